$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 was "Wait" (0 - 120) -> becomes "Idle01" (0 - 82)
$ws.Range("D4").Value = "Idle01"
$ws.Range("F4").Value = 82

# Row 5 was "Idle01~2" (121 - 291) -> split into "Idle02" (83 - 170)
$ws.Range("D5").Value = "Idle02"
$ws.Range("E5").Value = 83
$ws.Range("F5").Value = 170

# Row 6 was "Walk" (292 - 332) -> becomes "Wait" (171 - 291)
$ws.Range("D6").Value = "Wait"
$ws.Range("E6").Value = 171
$ws.Range("F6").Value = 291

# New row 7 holds "Walk" (292 - 332), matching the formatting of the rows above it
$ws.Range("C6:F6").Copy()
$ws.Range("C7:F7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Walk"
$ws.Range("E7").Value = 292
$ws.Range("F7").Value = 332

# Match the source workbook's final selection state
$ws.Range("F7").Select()
